$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 92
$ws.Range("H92").Value = 581.875
$ws.Range("I92").Value = 575.8333
$ws.Range("J92").Value = 600
$ws.Range("K92").Value = 575.8333
$ws.Range("L92").Value = 600
$ws.Range("M92").Value = 672.1667
$ws.Range("N92").Value = -3096

# Row 129
$ws.Range("H129").Value = 812.8570999999999
$ws.Range("I129").Value = 539.4
$ws.Range("J129").Value = 898.3125
$ws.Range("K129").Value = 1618.2
$ws.Range("L129").Value = 2694.9375
$ws.Range("M129").Value = 3381.8
$ws.Range("N129").Value = -12694.9375

# Row 137
$ws.Range("H137").Value = 16220.397
$ws.Range("I137").Value = 1382
$ws.Range("J137").Value = 38752.777
$ws.Range("K137").Value = 4146
$ws.Range("L137").Value = 116258.331
$ws.Range("M137").Value = -1596
$ws.Range("N137").Value = -121358.331

# Row 138
$ws.Range("H138").Value = 2527.966
$ws.Range("I138").Value = 2032.5
$ws.Range("J138").Value = 2682.111
$ws.Range("K138").Value = 6097.5
$ws.Range("L138").Value = 8046.333
$ws.Range("M138").Value = -957.5
$ws.Range("N138").Value = -18326.333

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 18644.064
$ws.Range("I32").Value = 19784.25
$ws.Range("J32").Value = 8002.3335
$ws.Range("K32").Value = 19784.25
$ws.Range("L32").Value = 8002.3335
$ws.Range("M32").Value = -19497.25
$ws.Range("N32").Value = -8576.333500000001

# Row 61
$ws.Range("H61").Value = 724726.0600000001
$ws.Range("I61").Value = 2255790.5
$ws.Range("J61").Value = 4225.1177
$ws.Range("K61").Value = 2255790.5
$ws.Range("L61").Value = 4225.1177
$ws.Range("M61").Value = -2255578.5
$ws.Range("N61").Value = -4649.1177

# Row 63
$ws.Range("H63").Value = 3474234.2
$ws.Range("I63").Value = 2263.625
$ws.Range("J63").Value = 31250000
$ws.Range("K63").Value = 2263.625
$ws.Range("L63").Value = 31250000
$ws.Range("M63").Value = -1577.625
$ws.Range("N63").Value = -31251372

# Row 66
$ws.Range("H66").Value = 3474234.2
$ws.Range("I66").Value = 2263.625
$ws.Range("J66").Value = 31250000
$ws.Range("K66").Value = 11318.125
$ws.Range("L66").Value = 156250000
$ws.Range("M66").Value = -7886.125
$ws.Range("N66").Value = -156256864

# Row 74
$ws.Range("H74").Value = 2304.5862
$ws.Range("I74").Value = 2704.1667
$ws.Range("J74").Value = 1650.7273
$ws.Range("K74").Value = 2704.1667
$ws.Range("L74").Value = 1650.7273
$ws.Range("M74").Value = -1830.1667
$ws.Range("N74").Value = -3398.7273

# Row 77
$ws.Range("H77").Value = 2304.5862
$ws.Range("I77").Value = 2704.1667
$ws.Range("J77").Value = 1650.7273
$ws.Range("K77").Value = 13520.8335
$ws.Range("L77").Value = 8253.636500000001
$ws.Range("M77").Value = -9152.833500000001
$ws.Range("N77").Value = -16989.6365

# Row 132
$ws.Range("H132").Value = 29128.21
$ws.Range("I132").Value = 2636
$ws.Range("J132").Value = 52971.2
$ws.Range("K132").Value = 7908
$ws.Range("L132").Value = 158913.6
$ws.Range("M132").Value = -5378
$ws.Range("N132").Value = -163973.6

# Row 135
$ws.Range("H135").Value = 19036.5
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 19036.5
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 19036.5
$ws.Range("N135").Value = -29176.5

# Row 136
$ws.Range("H136").Value = 724726.0600000001
$ws.Range("I136").Value = 2255790.5
$ws.Range("J136").Value = 4225.1177
$ws.Range("K136").Value = 6767371.5
$ws.Range("L136").Value = 12675.3531
$ws.Range("M136").Value = -6764821.5
$ws.Range("N136").Value = -17775.3531

$ws = $wb.Worksheets.Item("BSM")
# Row 64
$ws.Range("H64").Value = 963.3077
$ws.Range("I64").Value = 1680.5
$ws.Range("J64").Value = 348.57144
$ws.Range("K64").Value = 1680.5
$ws.Range("L64").Value = 348.57144
$ws.Range("M64").Value = -1455.5
$ws.Range("N64").Value = -798.5714399999999

# Row 67
$ws.Range("H67").Value = 963.3077
$ws.Range("I67").Value = 1680.5
$ws.Range("J67").Value = 348.57144
$ws.Range("K67").Value = 1680.5
$ws.Range("L67").Value = 348.57144
$ws.Range("M67").Value = -900.5
$ws.Range("N67").Value = -1908.57144

# Row 81
$ws.Range("H81").Value = 19273.143
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 19273.143
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 19273.143
$ws.Range("N81").Value = -21395.143

# Row 84
$ws.Range("H84").Value = 19273.143
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 19273.143
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 57819.429
$ws.Range("N84").Value = -68427.429

# Row 99
$ws.Range("H99").Value = 1755.4445
$ws.Range("I99").Value = 1787.375
$ws.Range("J99").Value = 1500
$ws.Range("K99").Value = 1787.375
$ws.Range("L99").Value = 1500
$ws.Range("M99").Value = -289.375
$ws.Range("N99").Value = -4496

# Row 134
$ws.Range("H134").Value = 38806.82
$ws.Range("I134").Value = 51375.76
$ws.Range("J134").Value = 1100
$ws.Range("K134").Value = 154127.28
$ws.Range("L134").Value = 3300
$ws.Range("M134").Value = -151592.28
$ws.Range("N134").Value = -8370

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 17437.715
$ws.Range("I31").Value = 34534.332
$ws.Range("J31").Value = 4615.25
$ws.Range("K31").Value = 34534.332
$ws.Range("L31").Value = 4615.25
$ws.Range("M31").Value = -34239.332
$ws.Range("N31").Value = -5205.25

# Row 34
$ws.Range("H34").Value = 17437.715
$ws.Range("I34").Value = 34534.332
$ws.Range("J34").Value = 4615.25
$ws.Range("K34").Value = 34534.332
$ws.Range("L34").Value = 4615.25
$ws.Range("M34").Value = -34332.332
$ws.Range("N34").Value = -5019.25

# Row 99
$ws.Range("H99").Value = 4749.2856
$ws.Range("I99").Value = 3380.0625
$ws.Range("J99").Value = 6574.9165
$ws.Range("K99").Value = 3380.0625
$ws.Range("L99").Value = 6574.9165
$ws.Range("M99").Value = -1882.0625
$ws.Range("N99").Value = -9570.916499999999

# Row 122
$ws.Range("H122").Value = 2383.7273
$ws.Range("I122").Value = 2535.6667
$ws.Range("J122").Value = 1700
$ws.Range("K122").Value = 7607.000100000001
$ws.Range("L122").Value = 5100
$ws.Range("M122").Value = -5157.000100000001
$ws.Range("N122").Value = -10000

# Row 126
$ws.Range("H126").Value = 4749.2856
$ws.Range("I126").Value = 3380.0625
$ws.Range("J126").Value = 6574.9165
$ws.Range("K126").Value = 10140.1875
$ws.Range("L126").Value = 19724.7495
$ws.Range("M126").Value = -7670.1875
$ws.Range("N126").Value = -24664.7495

$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 520
$ws.Range("I34").Value = 140
$ws.Range("J34").Value = 900
$ws.Range("K34").Value = 420
$ws.Range("L34").Value = 2700
$ws.Range("M34").Value = -336
$ws.Range("N34").Value = -2868

# Row 39
$ws.Range("H39").Value = 3101.3333
$ws.Range("I39").Value = 1900
$ws.Range("J39").Value = 3702
$ws.Range("K39").Value = 5700
$ws.Range("L39").Value = 11106
$ws.Range("M39").Value = -5406
$ws.Range("N39").Value = -11694

# Row 55
$ws.Range("H55").Value = 2375
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 2375
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 7125
$ws.Range("N55").Value = -7479

# Row 68
$ws.Range("H68").Value = 5114.44
$ws.Range("I68").Value = 873.5
$ws.Range("J68").Value = 5922.2383
$ws.Range("K68").Value = 2620.5
$ws.Range("L68").Value = 17766.7149
$ws.Range("M68").Value = -1809.5
$ws.Range("N68").Value = -19388.7149

# Row 71
$ws.Range("H71").Value = 5114.44
$ws.Range("I71").Value = 873.5
$ws.Range("J71").Value = 5922.2383
$ws.Range("K71").Value = 7861.5
$ws.Range("L71").Value = 53300.1447
$ws.Range("M71").Value = -3805.5
$ws.Range("N71").Value = -61412.1447

# Row 87
$ws.Range("H87").Value = 22900
$ws.Range("I87").Value = 14000
$ws.Range("J87").Value = 28833.334
$ws.Range("K87").Value = 42000
$ws.Range("L87").Value = 86500.00199999999
$ws.Range("M87").Value = -40752
$ws.Range("N87").Value = -88996.00199999999

# Row 90
$ws.Range("H90").Value = 22900
$ws.Range("I90").Value = 14000
$ws.Range("J90").Value = 28833.334
$ws.Range("K90").Value = 126000
$ws.Range("L90").Value = 259500.006
$ws.Range("M90").Value = -119760
$ws.Range("N90").Value = -271980.006

# Row 119
$ws.Range("H119").Value = 3406.889
$ws.Range("I119").Value = 1000
$ws.Range("J119").Value = 4094.5715
$ws.Range("K119").Value = 3000
$ws.Range("L119").Value = 12283.7145
$ws.Range("M119").Value = 1838
$ws.Range("N119").Value = -21959.7145

# Row 131
$ws.Range("H131").Value = 808.5306399999999
$ws.Range("I131").Value = 413.33334
$ws.Range("J131").Value = 821.0105
$ws.Range("K131").Value = 1240.00002
$ws.Range("L131").Value = 2463.0315
$ws.Range("M131").Value = 3799.99998
$ws.Range("N131").Value = -12543.0315

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 2355.5908
$ws.Range("I122").Value = 2195
$ws.Range("J122").Value = 2901.6
$ws.Range("K122").Value = 6585
$ws.Range("L122").Value = 8704.799999999999
$ws.Range("M122").Value = -4135
$ws.Range("N122").Value = -13604.8

# Row 132
$ws.Range("H132").Value = 112828.43
$ws.Range("I132").Value = 150371.58
$ws.Range("J132").Value = 75285.28999999999
$ws.Range("K132").Value = 451114.74
$ws.Range("L132").Value = 225855.87
$ws.Range("M132").Value = -448584.74
$ws.Range("N132").Value = -230915.87

$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 2287.6428
$ws.Range("I93").Value = 2335.5833
$ws.Range("J93").Value = 2000
$ws.Range("K93").Value = 2335.5833
$ws.Range("L93").Value = 2000
$ws.Range("M93").Value = -1087.5833
$ws.Range("N93").Value = -4496

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 2692.6
$ws.Range("I81").Value = 2050
$ws.Range("J81").Value = 5263
$ws.Range("K81").Value = 4100
$ws.Range("L81").Value = 10526
$ws.Range("M81").Value = -3039
$ws.Range("N81").Value = -12648

# Row 84
$ws.Range("H84").Value = 2692.6
$ws.Range("I84").Value = 2050
$ws.Range("J84").Value = 5263
$ws.Range("K84").Value = 20500
$ws.Range("L84").Value = 52630
$ws.Range("M84").Value = -15196
$ws.Range("N84").Value = -63238

# Row 100
$ws.Range("H100").Value = 933.2
$ws.Range("I100").Value = 450.25
$ws.Range("J100").Value = 1485.1428
$ws.Range("K100").Value = 900.5
$ws.Range("L100").Value = 2970.2856
$ws.Range("M100").Value = -359.5
$ws.Range("N100").Value = -4052.2856

# Row 122
$ws.Range("H122").Value = 1998.9375
$ws.Range("I122").Value = 1816.5834
$ws.Range("J122").Value = 2546
$ws.Range("K122").Value = 5449.7502
$ws.Range("L122").Value = 7638
$ws.Range("M122").Value = -2999.7502
$ws.Range("N122").Value = -12538

# Row 124
$ws.Range("H124").Value = 29500
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 29500
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 29500
$ws.Range("N124").Value = -39320

# Row 126
$ws.Range("H126").Value = 1067
$ws.Range("I126").Value = 855.93335
$ws.Range("J126").Value = 2650
$ws.Range("K126").Value = 2567.80005
$ws.Range("L126").Value = 7950
$ws.Range("M126").Value = -97.80004999999983
$ws.Range("N126").Value = -12890

# Row 132
$ws.Range("H132").Value = 2214
$ws.Range("I132").Value = 2026.15
$ws.Range("J132").Value = 2527.0833
$ws.Range("K132").Value = 6078.450000000001
$ws.Range("L132").Value = 7581.249899999999
$ws.Range("M132").Value = -3548.450000000001
$ws.Range("N132").Value = -12641.2499

# Row 133
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
